$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell carrying the default (General) style, used to restore
# the style on cells we temporarily mark as Text so numeric-looking
# strings (prices such as "592.99") are not coerced into numbers.
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.850.48"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = "  +3.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.638.32"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = "  +6.69%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.99"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.52"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.617.40"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = "  +6.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.606"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  +0.98%  "
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.203"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  +2.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.605"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  +1.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.93"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  +2.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000287"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "695.04"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.228.33"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  +6.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.99"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  +3.42%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.652.52"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  +7.21%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "71.992.99"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "  +3.58%  "
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.42"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  +3.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.58"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  +2.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.934"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  +2.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.80"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  +7.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.87"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  +3.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "103.71"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.03"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  +2.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.86"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  +4.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.97"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = "  +2.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.09"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  +2.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.12"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  +3.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.31"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  +4.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.16"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  +15.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "581.92"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  +3.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.34"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.109"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  +2.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.47"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.655.02"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.144"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.81"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0768"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  +6.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.43"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  +4.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0462"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  +8.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.76"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  +3.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.349"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  +2.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.42"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  +1.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.84"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  +5.64%  "
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.44"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  +3.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.33"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  -0.84%  "
